# Listas sem duplicação de professores
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "-"
$ws.Range("C3").Value = "[-, 'MCT-3A-Lab. de eletroeletrônica', -]"
$ws.Range("C4").Value = "[-, 'MCT-3A-Lab. de eletroeletrônica', -]"
$ws.Range("E4").Value = "[-, -, 'MCT-2A-Sistemas digitais']"
$ws.Range("C6").Value = "-"
$ws.Range("E6").Value = "[-, -, 'MCT-2A-Sistemas digitais']"
$ws.Range("D8").Value = "-"
$ws.Range("B20").Value = "-"
$ws.Range("C21").Value = "[-, 'ELM-2NA-Lab. Circuitos Elétricos']"
